$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set numeric values, replacing the shared-string text cells
$ws.Range("B3").Value = 4343
$ws.Range("A5").Value = 432.42340000000002

# Apply number format (numFmtId 2 => "0.00") to the changed cells
$ws.Range("B3").NumberFormat = "0.00"
$ws.Range("A5").NumberFormat = "0.00"

# Update the active cell selection to A5
$ws.Range("A5").Select()
